$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.687.10'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.468.63'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.07'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.96'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("E10").Value = '  +8.72%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.77'
$ws.Range("E12").Value = '  +4.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  +3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.019.26'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.58'
$ws.Range("E16").Value = '  -3.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.457.13'
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.08'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.673.59'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '462.80'
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.44'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.27'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.80'
$ws.Range("E25").Value = '  +15.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.32'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.16'
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.01'
$ws.Range("E34").Value = '  -4.51%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.19'
$ws.Range("E36").Value = '  +8.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0490'
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.06'
$ws.Range("E39").Value = '  +3.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '149.59'
$ws.Range("E40").Value = '  +4.75%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.323'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.72'
$ws.Range("E42").Value = '  +5.54%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.42'
$ws.Range("E45").Value = '  +3.93%  '
$ws.Range("E46").Value = '  +2.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0568'
$ws.Range("E47").Value = '  +33.97%  '
$ws.Range("E48").Value = '  +10.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.36'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.25'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("E51").Value = '  -4.39%  '
